$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 13 with the new timesheet entry
$ws.Range("A13").Value = 45577
$ws.Range("B13").Value = 0.48541666666666666
$ws.Range("C13").Value = 0.57152777777777775

# Expand the weekly-summary formulas in row 4 to include the new row 13 data
$ws.Range("M4").Formula = "=SUM(D10:D13)"
$ws.Range("N4").Formula = "=SUM(G10:G13)"

# Add weekly money total for week 1 (mirrors N3/N4 pattern)
$ws.Range("N2").Formula = "=G2"

$ws.Range("J10").Select()
